$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.1878
$ws.Range("A9").Value = -20.43749999999997
$ws.Range("C11").Value = -14.0597
$ws.Range("A18").Value = -23.16620000000001
$ws.Range("A20").Value = -22.17630000000003
$ws.Range("D21").Value = -7.484200000000003
